# Applies scheduled market-data refresh updates to the Leve profit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1063
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2350.963
$ws.Range("I132").Value = 1270.25
$ws.Range("J132").Value = 10996.667
$ws.Range("K132").Value = 3810.75
$ws.Range("L132").Value = 32990.001
$ws.Range("M132").Value = -1280.75
$ws.Range("N132").Value = -38050.001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1104.0834
$ws.Range("I135").Value = 901.4545000000001
$ws.Range("K135").Value = 8113.0905
$ws.Range("M135").Value = -5578.0905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2936.8462
$ws.Range("I45").Value = 2027.8572
$ws.Range("J45").Value = 3997.3333
$ws.Range("K45").Value = 2027.8572
$ws.Range("L45").Value = 3997.3333
$ws.Range("M45").Value = -1650.8572
$ws.Range("N45").Value = -4751.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1460.0322
$ws.Range("I110").Value = 1108.7059
$ws.Range("J110").Value = 1886.6428
$ws.Range("K110").Value = 1108.7059
$ws.Range("L110").Value = 1886.6428
$ws.Range("M110").Value = 936.2941000000001
$ws.Range("N110").Value = -5976.6428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1994.5
$ws.Range("I122").Value = 1994.5
$ws.Range("K122").Value = 5983.5
$ws.Range("M122").Value = -3533.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1200.8
$ws.Range("J64").Value = 1126
$ws.Range("L64").Value = 1126
$ws.Range("N64").Value = -1576

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 1200.8
$ws.Range("J67").Value = 1126
$ws.Range("L67").Value = 1126
$ws.Range("N67").Value = -2686

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1163.6666
$ws.Range("I86").Value = 990
$ws.Range("J86").Value = 1198.4
$ws.Range("K86").Value = 990
$ws.Range("L86").Value = 1198.4
$ws.Range("M86").Value = 133
$ws.Range("N86").Value = -3444.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1163.6666
$ws.Range("I89").Value = 990
$ws.Range("J89").Value = 1198.4
$ws.Range("K89").Value = 4950
$ws.Range("L89").Value = 5992
$ws.Range("M89").Value = 666
$ws.Range("N89").Value = -17224

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1247.3334
$ws.Range("I134").Value = 1247.3334
$ws.Range("K134").Value = 3742.0002
$ws.Range("M134").Value = -1207.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2166.3333
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2166.3333
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 2166.3333
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -2392.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 12000
$ws.Range("J106").Value = 12000
$ws.Range("L106").Value = 12000
$ws.Range("N106").Value = -14524

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2999
$ws.Range("I122").Value = 2999
$ws.Range("K122").Value = 8997
$ws.Range("M122").Value = -6547

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 123.8
$ws.Range("I38").Value = 100
$ws.Range("J38").Value = 159.5
$ws.Range("K38").Value = 300
$ws.Range("L38").Value = 478.5
$ws.Range("M38").Value = 47
$ws.Range("N38").Value = -1172.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 943.5
$ws.Range("J92").Value = 1000
$ws.Range("L92").Value = 3000
$ws.Range("N92").Value = -5496

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1407.0197
$ws.Range("I131").Value = 1536
$ws.Range("K131").Value = 4608
$ws.Range("M131").Value = 432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1869.8572
$ws.Range("I132").Value = 1869.8572
$ws.Range("K132").Value = 16828.7148
$ws.Range("M132").Value = -14298.7148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 19990
$ws.Range("I40").Value = 19990
$ws.Range("K40").Value = 19990
$ws.Range("M40").Value = -19839

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3999.25
$ws.Range("I80").Value = 3999
$ws.Range("K80").Value = 3999
$ws.Range("M80").Value = -3001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3999.25
$ws.Range("I83").Value = 3999
$ws.Range("K83").Value = 19995
$ws.Range("M83").Value = -15003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1654.4286
$ws.Range("I97").Value = 1071.875
$ws.Range("J97").Value = 2431.1667
$ws.Range("K97").Value = 1071.875
$ws.Range("L97").Value = 2431.1667
$ws.Range("M97").Value = -575.875
$ws.Range("N97").Value = -3423.1667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1152.7693
$ws.Range("I102").Value = 998.7
$ws.Range("K102").Value = 998.7
$ws.Range("M102").Value = 623.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 175572
$ws.Range("I23").Value = 194834
$ws.Range("J23").Value = 60000
$ws.Range("K23").Value = 194834
$ws.Range("L23").Value = 60000
$ws.Range("M23").Value = -194604
$ws.Range("N23").Value = -60460

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 9625
$ws.Range("I25").Value = 8125
$ws.Range("J25").Value = 11125
$ws.Range("K25").Value = 8125
$ws.Range("L25").Value = 11125
$ws.Range("M25").Value = -7895
$ws.Range("N25").Value = -11585

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3367.2354
$ws.Range("I46").Value = 2354.1667
$ws.Range("J46").Value = 5798.6
$ws.Range("K46").Value = 2354.1667
$ws.Range("L46").Value = 5798.6
$ws.Range("M46").Value = -2166.1667
$ws.Range("N46").Value = -6174.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 536.7222
$ws.Range("I55").Value = 494.3846
$ws.Range("J55").Value = 646.8
$ws.Range("K55").Value = 494.3846
$ws.Range("L55").Value = 646.8
$ws.Range("M55").Value = -321.3846
$ws.Range("N55").Value = -992.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1296.6666
$ws.Range("I93").Value = 1296.6666
$ws.Range("K93").Value = 1296.6666
$ws.Range("M93").Value = -48.66660000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2269.5
$ws.Range("I100").Value = 2274.5
$ws.Range("K100").Value = 2274.5
$ws.Range("M100").Value = -1733.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6388.684
$ws.Range("I122").Value = 4648.6665
$ws.Range("K122").Value = 13945.9995
$ws.Range("M122").Value = -11495.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 34642.715
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 65330
$ws.Range("J64").Value = 66500
$ws.Range("L64").Value = 66500
$ws.Range("N64").Value = -66996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 65330
$ws.Range("J67").Value = 66500
$ws.Range("L67").Value = 66500
$ws.Range("N67").Value = -68216

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1937.25
$ws.Range("I81").Value = 1937.25
$ws.Range("K81").Value = 3874.5
$ws.Range("M81").Value = -2813.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1937.25
$ws.Range("I84").Value = 1937.25
$ws.Range("K84").Value = 19372.5
$ws.Range("M84").Value = -14068.5
